$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts E..Z to F..AA)
$ws.Columns("E:E").Insert()

# New column header
$ws.Cells.Item(1, 5).Value2 = "areaCode"

# New column data (areaCode values for rows 2-7)
$ws.Cells.Item(2, 5).Value2 = 13
$ws.Cells.Item(3, 5).Value2 = 12
$ws.Cells.Item(4, 5).Value2 = 16
$ws.Cells.Item(5, 5).Value2 = 18
$ws.Cells.Item(6, 5).Value2 = 19
$ws.Cells.Item(7, 5).Value2 = 20

# Update the (now shifted) Location column F with row-specific values
$ws.Cells.Item(3, 6).Value2 = "12 11 13 11"
$ws.Cells.Item(4, 6).Value2 = "16 11 13 11"
$ws.Cells.Item(5, 6).Value2 = "18 11 13 11"
$ws.Cells.Item(6, 6).Value2 = "19 11 13 11"
$ws.Cells.Item(7, 6).Value2 = "20 11 13 11"

# Update the (now shifted) Status_Result column Z for a few rows
$ws.Cells.Item(3, 26).Value2 = 1
$ws.Cells.Item(4, 26).Value2 = 2
$ws.Cells.Item(6, 26).Value2 = 2

# Update view state to match the edited workbook
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W14").Select()
